# Applies metrics_4_5.xlsx edit: reorders the model names in column A
# (rows 2-26) and overwrites the metric columns B:Q with the same
# (new) constant set of values for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of model names for rows 2..26 (row -> model name)
$newNames = @{
    2  = "model_4_5_0"
    3  = "model_4_5_22"
    4  = "model_4_5_21"
    5  = "model_4_5_20"
    6  = "model_4_5_19"
    7  = "model_4_5_18"
    8  = "model_4_5_17"
    9  = "model_4_5_16"
    10 = "model_4_5_15"
    11 = "model_4_5_14"
    12 = "model_4_5_13"
    13 = "model_4_5_23"
    14 = "model_4_5_12"
    15 = "model_4_5_10"
    16 = "model_4_5_9"
    17 = "model_4_5_8"
    18 = "model_4_5_7"
    19 = "model_4_5_6"
    20 = "model_4_5_5"
    21 = "model_4_5_4"
    22 = "model_4_5_3"
    23 = "model_4_5_2"
    24 = "model_4_5_1"
    25 = "model_4_5_11"
    26 = "model_4_5_24"
}

# New constant values for columns B..Q (same for every data row)
$newValues = @{
    "B" = 0.445980542904741
    "C" = -0.05067851974480586
    "D" = -0.09972987042821879
    "E" = -2.197770350524112
    "F" = -0.07182324804279538
    "G" = 0.3288898821261522
    "H" = 0.6237281562693088
    "I" = 0.5329068068308546
    "J" = 0.3736368866304284
    "K" = 0.4532718467306415
    "L" = 0.279043433096329
    "M" = 0.5734892170966707
    "N" = 0.05025235926527027
    "O" = 0.5979038403548141
    "P" = 22.22406457805329
    "Q" = 34.41282282673529
}

for ($row = 2; $row -le 26; $row++) {
    $ws.Range("A$row").Value = $newNames[$row]
    foreach ($col in $newValues.Keys) {
        $ws.Range("$col$row").Value = $newValues[$col]
    }
}
